$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.89761366666667
$ws.Range("H2").Value = 38.692841
$ws.Range("I2").Value = 0.1295258291743358
$ws.Range("J2").Value = 0.1295258291743358
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 141.0704656666667
$ws.Range("N2").Value = 423.211397
$ws.Range("O2").Value = 0.05525296614535039
$ws.Range("P2").Value = 0.05525296614535039
$ws.Range("Q2").Value = 1819.472365945431
$ws.Range("R2").Value = 16375.25129350888
$ws.Range("S2").Value = 0.007156686254318014
$ws.Range("T2").Value = 0.007156686254318014

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.89761366666667
$ws.Range("H3").Value = 38.692841
$ws.Range("I3").Value = 0.1295258291743358
$ws.Range("J3").Value = 0.1295258291743358
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.309554666666667
$ws.Range("N3").Value = 9.928663999999999
$ws.Range("O3").Value = 0.001296250856544298
$ws.Range("P3").Value = 0.001296250856544298
$ws.Range("Q3").Value = 42.68535749938044
$ws.Range("R3").Value = 384.168217494424
$ws.Range("S3").Value = 0.0001678979670118433
$ws.Range("T3").Value = 0.0001678979670118433

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.89761366666667
$ws.Range("H4").Value = 38.692841
$ws.Range("I4").Value = 0.1295258291743358
$ws.Range("J4").Value = 0.1295258291743358
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2408.188354666666
$ws.Range("N4").Value = 7224.565063999999
$ws.Range("O4").Value = 0.9432133721485603
$ws.Range("P4").Value = 0.9432133721485604
$ws.Range("Q4").Value = 31059.88303505631
$ws.Range("R4").Value = 279538.9473155068
$ws.Range("S4").Value = 0.1221704941158636
$ws.Range("T4").Value = 0.1221704941158637

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.89761366666667
$ws.Range("H5").Value = 38.692841
$ws.Range("I5").Value = 0.1295258291743358
$ws.Range("J5").Value = 0.1295258291743358
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6061513333333334
$ws.Range("N5").Value = 1.818454
$ws.Range("O5").Value = 0.0002374108495449545
$ws.Range("P5").Value = 0.0002374108495449545
$ws.Range("Q5").Value = 7.817905720868223
$ws.Range("R5").Value = 70.361151487814
$ws.Range("S5").Value = [double]"3.075083714229371E-05"
$ws.Range("T5").Value = [double]"3.075083714229371E-05"

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 17.69923533333333
$ws.Range("H6").Value = 53.097706
$ws.Range("I6").Value = 0.1777466895466555
$ws.Range("J6").Value = 0.1777466895466556
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 141.0704656666667
$ws.Range("N6").Value = 423.211397
$ws.Range("O6").Value = 0.05525296614535039
$ws.Range("P6").Value = 0.05525296614535039
$ws.Range("Q6").Value = 2496.839370417254
$ws.Range("R6").Value = 22471.55433375528
$ws.Range("S6").Value = 0.009821031819969465
$ws.Range("T6").Value = 0.009821031819969466

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 17.69923533333333
$ws.Range("H7").Value = 53.097706
$ws.Range("I7").Value = 0.1777466895466555
$ws.Range("J7").Value = 0.1777466895466556
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.309554666666667
$ws.Range("N7").Value = 9.928663999999999
$ws.Range("O7").Value = 0.001296250856544298
$ws.Range("P7").Value = 0.001296250856544298
$ws.Range("Q7").Value = 58.57658689386489
$ws.Range("R7").Value = 527.189282044784
$ws.Range("S7").Value = 0.0002304042985727657
$ws.Range("T7").Value = 0.0002304042985727658

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 17.69923533333333
$ws.Range("H8").Value = 53.097706
$ws.Range("I8").Value = 0.1777466895466555
$ws.Range("J8").Value = 0.1777466895466556
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2408.188354666666
$ws.Range("N8").Value = 7224.565063999999
$ws.Range("O8").Value = 0.9432133721485603
$ws.Range("P8").Value = 0.9432133721485604
$ws.Range("Q8").Value = 42623.09241623813
$ws.Range("R8").Value = 383607.8317461432
$ws.Range("S8").Value = 0.1676530544355442
$ws.Range("T8").Value = 0.1676530544355443

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 17.69923533333333
$ws.Range("H9").Value = 53.097706
$ws.Range("I9").Value = 0.1777466895466555
$ws.Range("J9").Value = 0.1777466895466556
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6061513333333334
$ws.Range("N9").Value = 1.818454
$ws.Range("O9").Value = 0.0002374108495449545
$ws.Range("P9").Value = 0.0002374108495449545
$ws.Range("Q9").Value = 10.72841509628044
$ws.Range("R9").Value = 96.555735866524
$ws.Range("S9").Value = [double]"4.219899256907477E-05"
$ws.Range("T9").Value = [double]"4.219899256907477E-05"

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 30.16920766666667
$ws.Range("H10").Value = 90.507623
$ws.Range("I10").Value = 0.302977879439589
$ws.Range("J10").Value = 0.302977879439589
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 141.0704656666667
$ws.Range("N10").Value = 423.211397
$ws.Range("O10").Value = 0.05525296614535039
$ws.Range("P10").Value = 0.05525296614535039
$ws.Range("Q10").Value = 4255.984174331037
$ws.Range("R10").Value = 38303.85756897933
$ws.Range("S10").Value = 0.01674042651546566
$ws.Range("T10").Value = 0.01674042651546566

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 30.16920766666667
$ws.Range("H11").Value = 90.507623
$ws.Range("I11").Value = 0.302977879439589
$ws.Range("J11").Value = 0.302977879439589
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.309554666666667
$ws.Range("N11").Value = 9.928663999999999
$ws.Range("O11").Value = 0.001296250856544298
$ws.Range("P11").Value = 0.001296250856544298
$ws.Range("Q11").Value = 99.84664202285244
$ws.Range("R11").Value = 898.6197782056719
$ws.Range("S11").Value = 0.0003927353357375424
$ws.Range("T11").Value = 0.0003927353357375424

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 30.16920766666667
$ws.Range("H12").Value = 90.507623
$ws.Range("I12").Value = 0.302977879439589
$ws.Range("J12").Value = 0.302977879439589
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2408.188354666666
$ws.Range("N12").Value = 7224.565063999999
$ws.Range("O12").Value = 0.9432133721485603
$ws.Range("P12").Value = 0.9432133721485604
$ws.Range("Q12").Value = 72653.13457238697
$ws.Range("R12").Value = 653878.2111514828
$ws.Range("S12").Value = 0.2857727873526347
$ws.Range("T12").Value = 0.2857727873526347

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 30.16920766666667
$ws.Range("H13").Value = 90.507623
$ws.Range("I13").Value = 0.302977879439589
$ws.Range("J13").Value = 0.302977879439589
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6061513333333334
$ws.Range("N13").Value = 1.818454
$ws.Range("O13").Value = 0.0002374108495449545
$ws.Range("P13").Value = 0.0002374108495449545
$ws.Range("Q13").Value = 18.28710545276022
$ws.Range("R13").Value = 164.583949074842
$ws.Range("S13").Value = [double]"7.193023575108162E-05"
$ws.Range("T13").Value = [double]"7.19302357510816E-05"

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 38.809555
$ws.Range("H14").Value = 116.428665
$ws.Range("I14").Value = 0.3897496018394196
$ws.Range("J14").Value = 0.3897496018394196
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 141.0704656666667
$ws.Range("N14").Value = 423.211397
$ws.Range("O14").Value = 0.05525296614535039
$ws.Range("P14").Value = 0.05525296614535039
$ws.Range("Q14").Value = 5474.881996166113
$ws.Range("R14").Value = 49273.93796549501
$ws.Range("S14").Value = 0.02153482155559724
$ws.Range("T14").Value = 0.02153482155559725

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 38.809555
$ws.Range("H15").Value = 116.428665
$ws.Range("I15").Value = 0.3897496018394196
$ws.Range("J15").Value = 0.3897496018394196
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.309554666666667
$ws.Range("N15").Value = 9.928663999999999
$ws.Range("O15").Value = 0.001296250856544298
$ws.Range("P15").Value = 0.001296250856544298
$ws.Range("Q15").Value = 128.4423438615067
$ws.Range("R15").Value = 1155.98109475356
$ws.Range("S15").Value = 0.0005052132552221469
$ws.Range("T15").Value = 0.000505213255222147

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 38.809555
$ws.Range("H16").Value = 116.428665
$ws.Range("I16").Value = 0.3897496018394196
$ws.Range("J16").Value = 0.3897496018394196
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2408.188354666666
$ws.Range("N16").Value = 7224.565063999999
$ws.Range("O16").Value = 0.9432133721485603
$ws.Range("P16").Value = 0.9432133721485604
$ws.Range("Q16").Value = 93460.71840079551
$ws.Range("R16").Value = 841146.4656071595
$ws.Range("S16").Value = 0.3676170362445176
$ws.Range("T16").Value = 0.3676170362445177

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 38.809555
$ws.Range("H17").Value = 116.428665
$ws.Range("I17").Value = 0.3897496018394196
$ws.Range("J17").Value = 0.3897496018394196
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6061513333333334
$ws.Range("N17").Value = 1.818454
$ws.Range("O17").Value = 0.0002374108495449545
$ws.Range("P17").Value = 0.0002374108495449545
$ws.Range("Q17").Value = 23.52446350932334
$ws.Range("R17").Value = 211.72017158391
$ws.Range("S17").Value = [double]"9.253078408250435E-05"
$ws.Range("T17").Value = [double]"9.253078408250436E-05"

